$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("D3").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("F6").Value = "[-, -, 'MCT-3A-CAM', -]"
$ws.Range("D7").Value = "-"
$ws.Range("B8").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("C8").Value = "[-, -, 'MCT-3A-CAM', -]"
